# Rename the inline logo pictures that live in the document's headers and
# footers. The Pearson logo (in both footers) goes from "image1.png" to
# "image2.png", and the BTEC logo (in both headers) goes from "image2.jpg"
# to "image1.jpg".
#
# InlineShape objects don't live in Document.InlineShapes when they are
# inside a header/footer story - they have to be reached through the
# Section's HeaderFooter.Range.InlineShapes collection instead.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # Headers: primary (1) and first-page (2) both carry the BTEC logo.
    foreach ($idx in 1, 2, 3) {
        $hf = $sec.Headers.Item($idx)
        if ($hf.Exists) {
            for ($i = 1; $i -le $hf.Range.InlineShapes.Count; $i++) {
                $shp = $hf.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    # Footers: primary (1) and first-page (2) both carry the Pearson logo.
    foreach ($idx in 1, 2, 3) {
        $ft = $sec.Footers.Item($idx)
        if ($ft.Exists) {
            for ($i = 1; $i -le $ft.Range.InlineShapes.Count; $i++) {
                $shp = $ft.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
